$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 461.6
$ws.Range("I2").Value = 470.33334
$ws.Range("J2").Value = 448.5
$ws.Range("K2").Value = 470.33334
$ws.Range("L2").Value = 448.5
$ws.Range("M2").Value = -357.33334
$ws.Range("N2").Value = -674.5

$ws.Range("H17").Value = 2836.7693
$ws.Range("I17").Value = 3813
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 11439
$ws.Range("L17").Value = 6000
$ws.Range("M17").Value = -11271
$ws.Range("N17").Value = -6336

$ws.Range("H33").Value = 213.5
$ws.Range("J33").Value = 509.25
$ws.Range("L33").Value = 509.25
$ws.Range("N33").Value = -967.25

$ws.Range("H76").Value = 5923.864
$ws.Range("I76").Value = 4777.5
$ws.Range("J76").Value = 7299.5
$ws.Range("K76").Value = 4777.5
$ws.Range("L76").Value = 7299.5
$ws.Range("M76").Value = -4462.5
$ws.Range("N76").Value = -7929.5

$ws.Range("H79").Value = 5923.864
$ws.Range("I79").Value = 4777.5
$ws.Range("J79").Value = 7299.5
$ws.Range("K79").Value = 4777.5
$ws.Range("L79").Value = 7299.5
$ws.Range("M79").Value = -3685.5
$ws.Range("N79").Value = -9483.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 560.3333
$ws.Range("J4").Value = 797.5
$ws.Range("L4").Value = 797.5
$ws.Range("N4").Value = -1029.5

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws.Range("H61").Value = 1868.1818
$ws.Range("I61").Value = 1882.625
$ws.Range("J61").Value = 1829.6666
$ws.Range("K61").Value = 1882.625
$ws.Range("L61").Value = 1829.6666
$ws.Range("M61").Value = -1670.625
$ws.Range("N61").Value = -2253.6666

$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws.Range("H132").Value = 2583.9285
$ws.Range("I132").Value = 2379.5454
$ws.Range("J132").Value = 3333.3333
$ws.Range("K132").Value = 7138.6362
$ws.Range("L132").Value = 9999.999899999999
$ws.Range("M132").Value = -4608.6362
$ws.Range("N132").Value = -15059.9999

$ws.Range("H136").Value = 1868.1818
$ws.Range("I136").Value = 1882.625
$ws.Range("J136").Value = 1829.6666
$ws.Range("K136").Value = 5647.875
$ws.Range("L136").Value = 5488.9998
$ws.Range("M136").Value = -3097.875
$ws.Range("N136").Value = -10588.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 539.9167
$ws.Range("I80").Value = 461.375
$ws.Range("K80").Value = 461.375
$ws.Range("M80").Value = 536.625

$ws.Range("H83").Value = 539.9167
$ws.Range("I83").Value = 461.375
$ws.Range("K83").Value = 2306.875
$ws.Range("M83").Value = 2685.125

$ws.Range("H86").Value = 3795.6365
$ws.Range("I86").Value = 3078.6667
$ws.Range("K86").Value = 3078.6667
$ws.Range("M86").Value = -1955.6667

$ws.Range("H89").Value = 3795.6365
$ws.Range("I89").Value = 3078.6667
$ws.Range("K89").Value = 15393.3335
$ws.Range("M89").Value = -9777.333500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 3000
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H7").Value = 278.4
$ws.Range("I7").Value = 278.4
$ws.Range("K7").Value = 278.4
$ws.Range("M7").Value = -165.4

$ws.Range("H16").Value = 1138
$ws.Range("I16").Value = 973.6667
$ws.Range("K16").Value = 973.6667
$ws.Range("M16").Value = -686.6667

$ws.Range("H22").Value = 833.3333
$ws.Range("I22").Value = 450
$ws.Range("K22").Value = 450
$ws.Range("M22").Value = -100

$ws.Range("H62").Value = 69166.5
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 82399.8
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 82399.8
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -83647.8

$ws.Range("H65").Value = 69166.5
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 82399.8
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 411999
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -418239

$ws.Range("H113").Value = 1138
$ws.Range("I113").Value = 973.6667
$ws.Range("K113").Value = 973.6667
$ws.Range("M113").Value = 1196.3333

$ws.Range("H134").Value = 2191.7778
$ws.Range("I134").Value = 1288.3
$ws.Range("K134").Value = 3864.9
$ws.Range("M134").Value = -1329.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2790.75
$ws.Range("I131").Value = 1403.3334
$ws.Range("K131").Value = 4210.0002
$ws.Range("M131").Value = 829.9997999999996

$ws.Range("H132").Value = 1488
$ws.Range("I132").Value = 1467.1666
$ws.Range("K132").Value = 13204.4994
$ws.Range("M132").Value = -10674.4994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 50130.5
$ws.Range("J96").Value = 50130.5
$ws.Range("L96").Value = 50130.5
$ws.Range("N96").Value = -55622.5

$ws.Range("H97").Value = 1979.1111
$ws.Range("J97").Value = 1712.6
$ws.Range("L97").Value = 1712.6
$ws.Range("N97").Value = -2704.6

$ws.Range("H101").Value = 21972
$ws.Range("J101").Value = 21972
$ws.Range("L101").Value = 21972
$ws.Range("N101").Value = -28462

$ws.Range("H132").Value = 3145.182
$ws.Range("I132").Value = 2408.4285
$ws.Range("J132").Value = 4434.5
$ws.Range("K132").Value = 7225.2855
$ws.Range("L132").Value = 13303.5
$ws.Range("M132").Value = -4695.2855
$ws.Range("N132").Value = -18363.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2466.6667
$ws.Range("I22").Value = 1200
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 1200
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -905
$ws.Range("N22").Value = -5590

$ws.Range("H27").Value = 2466.6667
$ws.Range("I27").Value = 1200
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 1200
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -1093
$ws.Range("N27").Value = -5214

$ws.Range("H97").Value = 33333.332
$ws.Range("J97").Value = 33333.332
$ws.Range("L97").Value = 33333.332
$ws.Range("N97").Value = -35315.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 93300.5
$ws.Range("J46").Value = 93300.5
$ws.Range("L46").Value = 93300.5
$ws.Range("N46").Value = -93762.5

$ws.Range("H60").Value = 100000
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H113").Value = 1165.0667
$ws.Range("I113").Value = 1272.1428
$ws.Range("J113").Value = 1071.375
$ws.Range("K113").Value = 3816.4284
$ws.Range("L113").Value = 3214.125
$ws.Range("M113").Value = -1646.4284
$ws.Range("N113").Value = -7554.125

$ws.Range("H134").Value = 93300.5
$ws.Range("J134").Value = 93300.5
$ws.Range("L134").Value = 279901.5
$ws.Range("N134").Value = -284971.5

Write-Host "All changes applied"
